# ePhone_prices.xlsx — Excel-VBA assessment edit
# 1. Fill in the missing "ePhone 4s mini" price on the price-list sheet.
# 2. Add a new "Sheet1" worksheet (after the price list) that computes the
#    Mean and Median of the price column via live formulas.

$wb = $excel.ActiveWorkbook

# --- 1. ePhone Price List: fill in the last (previously blank) price ---
$priceList = $wb.Worksheets.Item(1)
$priceList.Range("B7").Value = 850

# --- 2. New "Sheet1" with summary statistics, placed after the price list ---
$summary = $wb.Worksheets.Add($null, $priceList)
$summary.Name = "Sheet1"

$summary.Range("A2").Value = "Mean"
$summary.Range("B2").Formula = "=AVERAGE('ePhone Price List'!B:B)"

$summary.Range("A3").Value = "Median"
$summary.Range("B3").Formula = "=MEDIAN('ePhone Price List'!B:B)"

# --- Leave the same cell selections behind that the author's session had ---
$priceList.Range("B8").Select() | Out-Null
$summary.Range("B4").Select() | Out-Null
